$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): F4 899 -> 901, F6 41 -> 42
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 901
$ws1.Range("F6").Value = 42

# Sheet "全部类型" (sheet4.xml): F5 899 -> 901, F7 41 -> 42
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 901
$ws4.Range("F7").Value = 42
